$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("topcis", $true, $false, $false, $false, $false, $true, 1, $false, "topics", 2)
